# Apply renamed sheets and updated file names (task orders)
$wb = $excel.ActiveWorkbook

# --- Rename sheets (worksheet names only; sheetId / order unchanged) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-1650996061083159"
$wb.Worksheets.Item(2).Name = "NB_TO-1650996062440149"
$wb.Worksheets.Item(3).Name = "RS_TO-1650996062440149"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509960624913397"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509960625793817"

# --- Sheet 1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509960610509908.csv"
$ws1.Range("B3").Value = "GNG_stims-1650996061067751.csv"
$ws1.Range("B4").Value = "go_stims-16509960610689428.csv"
$ws1.Range("B5").Value = "GNG_stims-1650996061081884.csv"

# --- Sheet 2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-1650996062146525.csv"
$ws2.Range("B3").Value = "ZB-match_8-16509960611798332.csv"
$ws2.Range("B4").Value = "TB-16509960624241474.csv"
$ws2.Range("B5").Value = "OB-16509960613337824.csv"
$ws2.Range("B6").Value = "ZB-match_9-16509960610898516.csv"
$ws2.Range("B7").Value = "ZB-match_3-16509960611267495.csv"
$ws2.Range("B8").Value = "TB-1650996061804661.csv"
$ws2.Range("B9").Value = "OB-16509960616194441.csv"
$ws2.Range("B10").Value = "OB-16509960617560792.csv"

# --- Sheet 3 (RS) --- (no cell value changes)

# --- Sheet 4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509960624591658.csv"
$ws4.Range("B3").Value = "ZM_stims-1650996062440149.csv"
$ws4.Range("B4").Value = "MM_stims-16509960624751606.csv"
$ws4.Range("B5").Value = "ZM_stims-16509960624591658.csv"
$ws4.Range("B6").Value = "MM_stims-16509960624913397.csv"
$ws4.Range("B7").Value = "ZM_stims-16509960624751606.csv"

# --- Sheet 5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16509960625331023.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509960625633826.csv"
$ws5.Range("B4").Value = "SAT_stims-1650996062492957.csv"
$ws5.Range("B5").Value = "SAT_stims-16509960625089948.csv"
